# Mark the "flowers" (row 43) and "hikes" (row 49) hub rows as done.
# Row 47 ("lighthouses") already has the exact target formatting pattern
# (DONE checked, feature/comment/original post marked "Y"), so we copy
# its cell formats across and then set the expected values/strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$templateRow = 47
$targetRows = @(43, 49)

foreach ($row in $targetRows) {
    # Column A: DONE checkbox -> "[X] "
    $ws.Range("A$templateRow").Copy() | Out-Null
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("A$row").Value = "[X] "

    # Column C: feature -> "Y"
    $ws.Range("C$templateRow").Copy() | Out-Null
    $ws.Range("C$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("C$row").Value = "Y"

    # Column D: comment -> "Y"
    $ws.Range("D$templateRow").Copy() | Out-Null
    $ws.Range("D$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("D$row").Value = "Y"

    # Column G: original post -> "Y"
    $ws.Range("G$templateRow").Copy() | Out-Null
    $ws.Range("G$row").PasteSpecial(-4122) | Out-Null
    $ws.Range("G$row").Value = "Y"
}

$excel.CutCopyMode = 0
